# Apply the cryptos-list refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.205.51"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "1.563.47"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'210.71"
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "'22.06"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "'0.249"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").Value = "'0.0596"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D12").Value = "1.783.91"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "1.551.99"
$ws.Range("E13").Value = "  -0.83%  "
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "27.179.31"
$ws.Range("E16").Value = "  +0.81%  "
$ws.Range("D17").Value = "'61.83"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "'7.44"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").Value = "'216.66"
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").Value = "0.0₃0702"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "'4.12"
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").Value = "'9.21"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").Value = "'1.94"
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("D25").Value = "'152.99"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").Value = "'6.63"
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'15.04"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "'0.107"
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("E30").Value = "  +2.39%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "'3.24"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("D33").Value = "'3.16"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").Value = "1.439.59"
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("D35").Value = "'1.11"
$ws.Range("E35").Value = "  +4.32%  "
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("D37").Value = "'2.34"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").Value = "'0.0166"
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("E40").Value = "  +1.45%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.40"
$ws.Range("E41").Value = "  +2.55%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'0.807"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'1.01"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "'0.998"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("D45").Value = "'64.23"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("D47").Value = "1.694.29"
$ws.Range("E47").Value = "  -0.19%  "
$ws.Range("D48").Value = "'85.46"
$ws.Range("E48").Value = "  -1.94%  "
$ws.Range("D49").Value = "'0.0525"
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("E51").Value = "  -0.99%  "
